{"js": "// Apply the Review_364 edit: swap in the \"Byte Latent Transformer\" review\n// content in place of the \"Large Concept Models\" review.\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst paras = body.paragraphs.items;\n\n// --- Paragraph 0: title line (two runs separated by a manual line break) ---\n// Change the date and the paper title, keeping the existing <w:br/> split.\nconst dateResults = body.search(\"19.12.24\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\ndateResults.items[0].insertText(\"18.12.24\", Word.InsertLocation.replace);\nawait context.sync();\n\nconst titleResults = body.search(\n  \"Large Concept Models: Language Modeling in a Sentence Representation Space\",\n  { matchCase: true }\n);\ntitleResults.load(\"items\");\nawait context.sync();\ntitleResults.items[0].insertText(\n  \"Byte Latent Transformer: Patches Scale Better Than Tokens\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- Paragraphs 1-8 (existing): replace the body text of each paragraph. ---\n// Original paragraph 8 was the arXiv link; its text becomes the new\n// \"\u05db\u05dc \u05d0\u05dc\u05d5 \u05de\u05d5\u05db\u05e0\u05e1\u05d9\u05dd...\" paragraph below, and brand-new paragraphs (including a\n// new link) are appended after it further down.\nconst newBodyTexts = [\n  \"\u05db\u05de\u05d5\u05d1\u05df \u05dc\u05d0 \u05d9\u05db\u05d5\u05dc\u05ea\u05d9 \u05dc\u05e4\u05e1\u05e4\u05e1 \u05d0\u05ea \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05e9\u05d4\u05ea\u05e4\u05e8\u05e1\u05dd \u05dc\u05e4\u05e0\u05d9 \u05db\u05de\u05d4 \u05d9\u05de\u05d9\u05dd \u05d5\u05d2\u05e8\u05dd \u05dc\u05dc\u05d0 \u05de\u05e2\u05d8 \u05ea\u05d4\u05d5\u05d3\u05d4 \u05d1\u05e7\u05d4\u05d9\u05dc\u05ea AI. \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05dc\u05d4\u05d7\u05dc\u05d9\u05e3 \u05d0\u05ea \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05d9\u05d6\u05e8 \u05d4\u05e1\u05d8\u05d8\u05d9 \u05e9\u05d9\u05e9 \u05d1\u05db\u05dc \u05de\u05d5\u05d3\u05dc \u05d4\u05e9\u05e4\u05d4 \u05d1\u05de\u05e0\u05d2\u05e0\u05d5\u05df \u05d3\u05d9\u05e0\u05d0\u05de\u05d9 \u05e9\u05d1\u05d5\u05e0\u05d4 \u05d0\u05ea \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05d7\u05d3\u05e9\u05d9\u05dd (\u05e9\u05e7\u05d9\u05d1\u05dc\u05d5 \u05e9\u05dd \u05e4\u05d0\u05e6'\u05d9\u05dd) \u05db\u05dc\u05d5\u05de\u05e8 \u05db\u05d6\u05d4 \u05e9\u05d1\u05d5\u05e0\u05d4 \u05d0\u05d5\u05ea\u05dd \u05d1\u05ea\u05dc\u05d5\u05ea \u05d1\u05d4\u05e7\u05e9\u05e8 (contextualized). \",\n  \"\u05d4\u05e8\u05e6\u05d9\u05d5\u05e0\u05dc \u05db\u05d0\u05df \u05d4\u05d5\u05d0 \u05d3\u05d9 \u05d1\u05e8\u05d5\u05e8 \u05d4\u05e8\u05d9 \u05dc\u05e4\u05e2\u05de\u05d9\u05dd \u05d9\u05e9 \u05de\u05e7\u05e8\u05d9\u05dd \u05e9\u05d7\u05d9\u05d6\u05d5\u05d9 \u05e9\u05dc \u05db\u05de\u05d4 \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05d1\u05d0\u05d9\u05dd \u05d4\u05d5\u05d0 \u05d3\u05d9 \u05d1\u05e8\u05d5\u05e8 \u05d5\u05e0\u05d9\u05ea\u05df \u05dc\u05e2\u05e9\u05d5\u05ea \u05d0\u05d5\u05ea\u05d4 \u05db\u05de\u05e7\u05e9\u05d4 \u05d0\u05d7\u05ea (\u05db\u05dc\u05d5\u05de\u05e8 \u05dc\u05d0\u05d7\u05d3 \u05d0\u05ea \u05db\u05dc \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05dc\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d0\u05d7\u05d3 \u05d0\u05e8\u05d5\u05da \u05d0\u05d5 \u05e4\u05d0\u05e5' \u05dc\u05e4\u05d9 \u05e9\u05de\u05d5 \u05d1\u05de\u05d0\u05de\u05e8). \u05d5\u05dc\u05e4\u05e2\u05de\u05d9\u05dd \u05d4\u05de\u05e6\u05d1 \u05d4\u05d5\u05d0 \u05d4\u05e4\u05d5\u05da \u05d5\u05d4\u05d9\u05d9\u05e0\u05d5 \u05e8\u05d5\u05e6\u05d9\u05dd \u05dc\u05d7\u05d6\u05d5\u05ea \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d1\u05d2\u05e8\u05e0\u05d5\u05dc\u05e8\u05d9\u05d5\u05ea \u05e7\u05d8\u05e0\u05d4 \u05d9\u05d5\u05ea\u05e8. \u05d5\u05db\u05de\u05d5\u05d1\u05df \u05e9\u05d6\u05d4 \u05d1\u05dc\u05ea\u05d9 \u05d0\u05e4\u05e9\u05e8\u05d9 \u05d1\u05de\u05d5\u05d3\u05dc \u05e9\u05d9\u05e9 \u05d1\u05d4\u05dd \u05de\u05d9\u05dc\u05d5\u05df \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e7\u05d1\u05d5\u05e2.\",\n  \"\u05db\u05d0\u05de\u05d5\u05e8 \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05dc\u05d4\u05db\u05e0\u05d9\u05e1 \u05d3\u05d9\u05e0\u05de\u05d9\u05d5\u05ea \u05d1\u05d1\u05e0\u05d9\u05d9\u05ea \u05e4\u05d0\u05e6'\u05d9\u05dd (\u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05d7\u05d3\u05e9\u05d9\u05dd). \u05d0\u05d9\u05da \u05d4\u05d5\u05d0 \u05e2\u05d5\u05e9\u05d4 \u05d0\u05ea \u05d6\u05d4. \u05dc\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e0\u05ea\u05d5\u05df \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc \u05e8\u05d3\u05d5\u05d3 \u05d9\u05d7\u05e1\u05d9\u05ea \u05d1\u05e8\u05de\u05d4 \u05e9\u05dc \u05d1\u05d8\u05d9\u05dd (bytes) \u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05d8\u05e8\u05d4 \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05d9\u05d0 \u05dc\u05d7\u05d6\u05d5\u05ea \u05d0\u05ea \u05d4\u05d1\u05d9\u05d9\u05d8 \u05d4\u05d1\u05d0. \u05d5\u05d0\u05d6 \u05d1\u05de\u05d5\u05d3\u05dc \u05d4\u05d2\u05d3\u05d5\u05dc \u05e9\u05dc\u05e0\u05d5 \u05d4\u05dd \u05e7\u05d5\u05d1\u05e2\u05d9\u05dd \u05d0\u05ea \u05d2\u05d1\u05d5\u05dc\u05d5\u05ea \u05d4\u05e4\u05d0\u05e5 \u05e2\u05dc \u05e1\u05de\u05da \u05d0\u05e0\u05d8\u05e8\u05d5\u05e4\u05d9\u05d4 \u05e9\u05dc \u05d4\u05d1\u05d8\u05d9\u05dd. \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05dd \u05d4\u05d0\u05e0\u05d8\u05e8\u05d5\u05e4\u05d9\u05d4 \u05e9\u05dc \u05d4\u05d1\u05d9\u05d9\u05d8 \u05d0\u05d5 \u05d2\u05d3\u05d5\u05dc\u05d4 \u05de\u05e1\u05e3 \u05de\u05e1\u05d5\u05d9\u05dd \u05d0\u05d5 \u05d7\u05d5\u05d5\u05ea\u05d4 \u05e2\u05dc\u05d9\u05d4 \u05de\u05e2\u05dc \u05e1\u05e3 \u05de\u05e1\u05d5\u05d9\u05dd \u05de\u05e2\u05dc \u05d4\u05d0\u05e0\u05d8\u05e8\u05d5\u05e4\u05d9\u05d4 \u05e9\u05dc \u05d4\u05d1\u05d9\u05d9\u05d8 \u05d4\u05d1\u05d0, \u05e4\u05d5\u05ea\u05d7\u05d9\u05dd \u05e4\u05d0\u05e5' \u05d7\u05d3\u05e9. \u05d0\u05d7\u05e8\u05ea \u05de\u05de\u05e9\u05d9\u05db\u05d9\u05dd \u05d0\u05ea \u05d4\u05e4\u05d0\u05e5' \u05d4\u05e0\u05d5\u05db\u05d7\u05d9.\",\n  \"\u05d0\u05d1\u05dc \u05d0\u05d9\u05da \u05db\u05dc \u05d4\u05e1\u05d9\u05e4\u05d5\u05e8 \u05d4\u05d6\u05d4 \u05e2\u05d5\u05d1\u05d3 - \u05db\u05de\u05d5 \u05e9\u05d0\u05de\u05e8\u05ea\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05d5\u05d0 byte-level \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05d5\u05d0 \u05de\u05d0\u05d5\u05de\u05df \u05dc\u05d7\u05d6\u05d5\u05ea \u05d0\u05ea \u05d4\u05d1\u05d9\u05d9\u05d8 \u05d4\u05d1\u05d0 \u05d1\u05d8\u05e7\u05e1\u05d8. \u05d0\u05d1\u05dc \u05d1\u05de\u05e7\u05d5\u05dd \u05dc\u05d4\u05e1\u05ea\u05db\u05dc \u05e2\u05dc \u05d4\u05e7\u05d5\u05e0\u05e7\u05e1\u05d8 \u05d1\u05ea\u05d5\u05e8 \u05de\u05e2\u05e8\u05da \u05e9\u05dc \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05dc\u05d4\u05d7\u05dc\u05d9\u05e3 \u05d0\u05d5\u05ea\u05d5 \u05d1\u05e4\u05d0\u05e6\u05d9\u05dd \u05d3\u05d9\u05e0\u05de\u05d9\u05d9\u05dd \u05e0\u05e7\u05d1\u05e2\u05d9\u05dd \u05e2\u05dc \u05e1\u05de\u05da \u05d4\u05d0\u05e0\u05d8\u05e8\u05d5\u05e4\u05d9\u05d4 \u05db\u05de\u05d5 \u05e9\u05d4\u05e1\u05d1\u05e8\u05ea\u05d9 \u05e7\u05d5\u05d3\u05dd. \",\n  \"\u05d1\u05e0\u05d5\u05e1\u05e3 \u05dc\u05e4\u05d0\u05e6\u05d9\u05dd \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e9\u05ea\u05de\u05e9 \u05d2\u05dd \u05d1\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05d1\u05d8\u05d9\u05dd \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea n-grams (\u05dc\u05d5\u05e7\u05d7\u05d9\u05dd n-grams \u05dc\u05d1\u05d9\u05d9\u05d8 \u05e0\u05ea\u05d5\u05df \u05de n=3 \u05e2\u05d3 n=8, \u05de\u05e4\u05e2\u05dc\u05d9\u05dd \u05d0\u05d9\u05d6\u05d4 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05d0\u05e9, \u05e1\u05d5\u05db\u05de\u05d9\u05dd \u05d5\u05de\u05e0\u05e8\u05de\u05dc\u05d9\u05dd). \u05d0\u05ea \u05d4\u05ea\u05d5\u05e6\u05d0\u05d4 \u05d4\u05d5\u05e4\u05db\u05d9\u05dd \u05dc\u05d5\u05d5\u05e7\u05d8\u05d5\u05e8 (\u05d4\u05de\u05d0\u05de\u05e8 \u05dc\u05d0 \u05de\u05e4\u05e8\u05e9 \u05d0\u05d9\u05da- \u05e8\u05e7 \u05de\u05d6\u05db\u05d9\u05e8 \u05e9\u05d9\u05e9 \u05d0\u05d9\u05d6\u05d5 \u05e9\u05db\u05d1\u05d4 \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea \u05d4\u05de\u05e2\u05d5\u05e8\u05d1\u05ea \u05d1\u05d6\u05d4) \u05d5\u05de\u05d6\u05d9\u05df \u05d0\u05d5\u05ea\u05d5 \u05dc\u05de\u05d4 \u05e9\u05e7\u05e8\u05d5\u05d9 \u05d1\u05de\u05d0\u05de\u05e8 Encoder Multi-Headed Cross-Attention (\u05e0\u05e7\u05e8\u05d0 \u05dc\u05d6\u05d4 \u05dc\u05e4\u05e9\u05d8\u05d5\u05ea EMHCA). \",\n  \"\u05de\u05d8\u05e8\u05ea\u05d5 \u05e9\u05dc EMHCA \u05d4\u05d9\u05d0 \u05dc\u05e9\u05dc\u05d1 \u05d0\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05e4\u05d0\u05e6'\u05d9\u05dd \u05e2\u05dd \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d1\u05d8\u05d9\u05dd \u05e9\u05dc\u05d4\u05dd(\u05db\u05dc \u05e4\u05d0\u05e5 \u05de\u05ea\u05d7\u05e9\u05d1 \u05e8\u05e7 \u05d1\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d1\u05d8\u05d9\u05dd \u05e9\u05dc\u05d5 \u05d5\u05dc\u05d0 \u05e9\u05dc \u05d4\u05d0\u05d7\u05e8\u05d9\u05dd). \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05d4\u05ea\u05d7\u05dc\u05ea\u05d9 \u05e9\u05dc \u05db\u05dc \u05e4\u05d0\u05e5 \u05de\u05d7\u05d5\u05e9\u05d1 \u05db-pooling (\u05db\u05dc\u05d5\u05de\u05e8 \u05de\u05de\u05d5\u05e6\u05e2) \u05e9\u05dc \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d1\u05d8\u05d9\u05dd \u05e9\u05dc\u05d5 (\u05e0\u05d6\u05db\u05d9\u05e8 \u05d6\u05d4 \u05db\u05dc \u05e4\u05d0\u05e5 \u05d4\u05d9\u05e0\u05d5 \u05de\u05e2\u05e8\u05da \u05e9\u05dc \u05d4\u05d1\u05d8\u05d9\u05dd). \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05e0\u05d5 \u05d1\u05d5\u05e0\u05d9\u05dd \u05db\u05db\u05d4 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05db\u05dc \u05e4\u05d0\u05e5' \u05d4\u05de\u05ea\u05d7\u05e9\u05d1 \u05e8\u05e7 \u05d1\u05de\u05d4 \u05e9\u05d9\u05e9 \u05d1\u05ea\u05d5\u05db\u05d5 (internal representation).  \",\n  \"\u05d0\u05d6 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05d1\u05d8\u05d9\u05dd \u05d5\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05e4\u05d0\u05e6'\u05d9\u05dd \u05de\u05d5\u05d6\u05e0\u05d9\u05dd \u05dc-EMHCA \u05e9\u05d6\u05d4 \u05dc\u05de\u05e2\u05e9\u05d4 \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05d3\u05d9 \u05e8\u05d3\u05d5\u05d3 (\u05e2\u05dd \u05de\u05e2\u05d8 \u05e9\u05db\u05d1\u05d5\u05ea) \u05e9\u05de\u05d8\u05e8\u05d5\u05ea \u05dc\u05d1\u05e0\u05d5\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2 \u05ea\u05dc\u05d5\u05d9 \u05d4\u05e7\u05e9\u05e8 \u05e9\u05e4\u05d0\u05e6'\u05d9\u05dd \u05db\u05ea\u05dc\u05d5\u05ea \u05d1\u05d1\u05d8\u05d9\u05dd \u05e9\u05dc\u05d5. \u05db\u05dc\u05d5\u05de\u05e8 \u05d2\u05dd \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d1\u05d8\u05d9\u05dd \u05d4\u05dd keys and values \u05db\u05d0\u05df \u05db\u05d0\u05e9\u05e8 \u05d4-queries \u05d4\u05dd \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05e4\u05d0\u05e6\u05d9\u05dd. \u05db\u05d0\u05de\u05d5\u05e8 \u05de\u05d4 \u05e9\u05d9\u05d5\u05e6\u05d0 \u05de\u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05d4\u05e8\u05d3\u05d5\u05d3 \u05d4\u05d6\u05d4 \u05d4\u05d5\u05d0 \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05e4\u05d0\u05e6'\u05d9\u05dd. \u05e0\u05e6\u05d9\u05d9\u05df \u05e9- EMHCA \u05e4\u05d5\u05dc\u05d8 \u05d2\u05dd \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d1\u05d9\u05d8\u05d9\u05dd \u05d1\u05e1\u05d5\u05e3 (\u05dc\u05d0 \u05d4\u05e6\u05dc\u05d7\u05ea\u05d9 \u05dc\u05d4\u05d1\u05d9\u05df \u05d0\u05d9\u05da \u05d6\u05d4 \u05e0\u05d1\u05e0\u05d4). \",\n  \"\u05db\u05dc \u05d0\u05dc\u05d5 \u05de\u05d5\u05db\u05e0\u05e1\u05d9\u05dd \u05dc\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05d9\u05d5\u05ea\u05e8 \u05e2\u05de\u05d5\u05e7 \u05d5\u05db\u05d1\u05d3 \u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05ea \u05d4\u05d9\u05d5\u05e6\u05e8 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d9\u05d5\u05ea\u05e8 \u05f4\u05e2\u05de\u05d5\u05e7\u05f4 \u05e9\u05dc \u05d4\u05e4\u05d0\u05e6\u05d9\u05dd. \u05d1\u05e9\u05dc\u05d1 \u05d4\u05d0\u05d7\u05e8\u05d5\u05df \u05d9\u05e9 \u05d0\u05ea \u05d4-Local Decoder \u05e9\u05d4\u05d5\u05e4\u05da \u05d0\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05e4\u05d0\u05e6'\u05d9\u05dd \u05d9\u05d7\u05d3 \u05e2\u05dd \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d1\u05d8\u05d9\u05dd \u05dc\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d1\u05d8\u05d9\u05dd \u05d4\u05e1\u05d5\u05e4\u05d9\u05d9\u05dd \u05e9\u05de\u05d4\u05dd \u05e0\u05d7\u05d6\u05d4 \u05d4\u05d1\u05d9\u05d9\u05d8 \u05d4\u05d1\u05d0. \u05d6\u05d4 \u05d2\u05dd \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05e8\u05d3\u05d5\u05d3 \u05d0\u05d1\u05dc \u05d4\u05e4\u05e2\u05dd \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05e4\u05d0\u05e6'\u05d9\u05dd \u05d4\u05dd keys and values \u05d5\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d1\u05d8\u05d9\u05dd \u05d4\u05dd \u05d4-queries.\",\n];\n\nfor (let i = 0; i < newBodyTexts.length; i++) {\n  paras[i + 1].insertText(newBodyTexts[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- Append two brand-new paragraphs, then a brand-new arXiv link paragraph,\n// after the paragraph that used to hold the (now-replaced) old link text. ---\nbody.paragraphs.load(\"items\");\nawait context.sync();\nconst refreshedParas = body.paragraphs.items;\nconst lastPara = refreshedParas[refreshedParas.length - 1];\n\nconst newPara1 = lastPara.insertParagraph(\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d8\u05d5\u05e2\u05df \u05dc\u05db\u05dc \u05de\u05d9\u05e0\u05d9 \u05d9\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05e9\u05dc \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05de\u05d5\u05e6\u05e2\u05ea \u05db\u05de\u05d5 \u05d9\u05db\u05d5\u05dc\u05ea \u05dc\u05d7\u05d6\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05dc\u05e2\u05dc\u05d5\u05ea \u05d0\u05d9\u05e0\u05e4\u05e8\u05e0\u05e1 \u05e7\u05d1\u05d5\u05e2\u05d4, \u05d5\u05de\u05e6\u05d9\u05d2\u05d4 \u05d3\u05d9\u05d5\u05e7 \u05de\u05e9\u05d5\u05e4\u05e8 \u05d1\u05d0\u05d9\u05de\u05d5\u05df \u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\nconst newPara2 = newPara1.insertParagraph(\n  \"\u05d0\u05d5\u05e7\u05d9\u05d9, \u05d7\u05d9\u05d9\u05d1 \u05dc\u05d4\u05d2\u05d9\u05d3 \u05e9\u05d4\u05de\u05d0\u05de\u05e8 \u05dc\u05d0 \u05db\u05ea\u05d5\u05d1 \u05db\u05d6\u05d4 \u05d8\u05d5\u05d1 - \u05d9\u05e9 \u05d3\u05d1\u05e8\u05d9\u05dd \u05e9\u05dc\u05d0 \u05d4\u05d5\u05e1\u05d1\u05e8\u05d5 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d1\u05e8\u05d5\u05e8\u05d4 (\u05dc\u05de\u05d9\u05d8\u05d1 \u05d9\u05d3\u05d9\u05e2\u05ea\u05d9 \u05db\u05de\u05d5\u05d1\u05df). \u05d0\u05e0\u05d9 \u05e8\u05e7 \u05de\u05e7\u05d5\u05d5\u05d4 \u05e9\u05d4\u05e6\u05dc\u05d7\u05ea\u05d9 \u05dc\u05d4\u05d1\u05d9\u05df \u05d0\u05d5\u05ea\u05d5 \u05e0\u05db\u05d5\u05df\u2026.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\nconst newPara3 = newPara2.insertParagraph(\n  \"https://arxiv.org/abs/2412.09871\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "# Apply the Review_364 edit: swap in the \"Byte Latent Transformer\" review\n# content in place of the \"Large Concept Models\" review.\n\n$d = $word.ActiveDocument\n\n# --- Paragraph 1: title line (two runs separated by a manual line break) ---\n# Change the date and the paper title via Find/Replace so the existing\n# <w:br/> line-break split between the two runs is left untouched.\n$find1 = $d.Content\n$find1.Find.Execute(\"19.12.24\", $false, $false, $false, $false, $false, $true, 1, $false, \"18.12.24\", 2)\n\n$find2 = $d.Content\n$find2.Find.Execute(\"Large Concept Models: Language Modeling in a Sentence Representation Space\", $false, $false, $false, $false, $false, $true, 1, $false, \"Byte Latent Transformer: Patches Scale Better Than Tokens\", 2)\n\n# --- Paragraphs 2-9 (existing): replace the body text of each paragraph. ---\n# Original paragraph 9 was the arXiv link; its text becomes the new\n# \"\u05db\u05dc \u05d0\u05dc\u05d5 \u05de\u05d5\u05db\u05e0\u05e1\u05d9\u05dd...\" paragraph below, and brand-new paragraphs (including a\n# new link) are appended after it further down.\n$d.Paragraphs.Item(2).Range.Text = \"\u05db\u05de\u05d5\u05d1\u05df \u05dc\u05d0 \u05d9\u05db\u05d5\u05dc\u05ea\u05d9 \u05dc\u05e4\u05e1\u05e4\u05e1 \u05d0\u05ea \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05e9\u05d4\u05ea\u05e4\u05e8\u05e1\u05dd \u05dc\u05e4\u05e0\u05d9 \u05db\u05de\u05d4 \u05d9\u05de\u05d9\u05dd \u05d5\u05d2\u05e8\u05dd \u05dc\u05dc\u05d0 \u05de\u05e2\u05d8 \u05ea\u05d4\u05d5\u05d3\u05d4 \u05d1\u05e7\u05d4\u05d9\u05dc\u05ea AI. \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05dc\u05d4\u05d7\u05dc\u05d9\u05e3 \u05d0\u05ea \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05d9\u05d6\u05e8 \u05d4\u05e1\u05d8\u05d8\u05d9 \u05e9\u05d9\u05e9 \u05d1\u05db\u05dc \u05de\u05d5\u05d3\u05dc \u05d4\u05e9\u05e4\u05d4 \u05d1\u05de\u05e0\u05d2\u05e0\u05d5\u05df \u05d3\u05d9\u05e0\u05d0\u05de\u05d9 \u05e9\u05d1\u05d5\u05e0\u05d4 \u05d0\u05ea \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05d7\u05d3\u05e9\u05d9\u05dd (\u05e9\u05e7\u05d9\u05d1\u05dc\u05d5 \u05e9\u05dd \u05e4\u05d0\u05e6'\u05d9\u05dd) \u05db\u05dc\u05d5\u05de\u05e8 \u05db\u05d6\u05d4 \u05e9\u05d1\u05d5\u05e0\u05d4 \u05d0\u05d5\u05ea\u05dd \u05d1\u05ea\u05dc\u05d5\u05ea \u05d1\u05d4\u05e7\u05e9\u05e8 (contextualized). \"\n$d.Paragraphs.Item(3).Range.Text = \"\u05d4\u05e8\u05e6\u05d9\u05d5\u05e0\u05dc \u05db\u05d0\u05df \u05d4\u05d5\u05d0 \u05d3\u05d9 \u05d1\u05e8\u05d5\u05e8 \u05d4\u05e8\u05d9 \u05dc\u05e4\u05e2\u05de\u05d9\u05dd \u05d9\u05e9 \u05de\u05e7\u05e8\u05d9\u05dd \u05e9\u05d7\u05d9\u05d6\u05d5\u05d9 \u05e9\u05dc \u05db\u05de\u05d4 \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05d1\u05d0\u05d9\u05dd \u05d4\u05d5\u05d0 \u05d3\u05d9 \u05d1\u05e8\u05d5\u05e8 \u05d5\u05e0\u05d9\u05ea\u05df \u05dc\u05e2\u05e9\u05d5\u05ea \u05d0\u05d5\u05ea\u05d4 \u05db\u05de\u05e7\u05e9\u05d4 \u05d0\u05d7\u05ea (\u05db\u05dc\u05d5\u05de\u05e8 \u05dc\u05d0\u05d7\u05d3 \u05d0\u05ea \u05db\u05dc \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05dc\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d0\u05d7\u05d3 \u05d0\u05e8\u05d5\u05da \u05d0\u05d5 \u05e4\u05d0\u05e5' \u05dc\u05e4\u05d9 \u05e9\u05de\u05d5 \u05d1\u05de\u05d0\u05de\u05e8). \u05d5\u05dc\u05e4\u05e2\u05de\u05d9\u05dd \u05d4\u05de\u05e6\u05d1 \u05d4\u05d5\u05d0 \u05d4\u05e4\u05d5\u05da \u05d5\u05d4\u05d9\u05d9\u05e0\u05d5 \u05e8\u05d5\u05e6\u05d9\u05dd \u05dc\u05d7\u05d6\u05d5\u05ea \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d1\u05d2\u05e8\u05e0\u05d5\u05dc\u05e8\u05d9\u05d5\u05ea \u05e7\u05d8\u05e0\u05d4 \u05d9\u05d5\u05ea\u05e8. \u05d5\u05db\u05de\u05d5\u05d1\u05df \u05e9\u05d6\u05d4 \u05d1\u05dc\u05ea\u05d9 \u05d0\u05e4\u05e9\u05e8\u05d9 \u05d1\u05de\u05d5\u05d3\u05dc \u05e9\u05d9\u05e9 \u05d1\u05d4\u05dd \u05de\u05d9\u05dc\u05d5\u05df \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e7\u05d1\u05d5\u05e2.\"\n$d.Paragraphs.Item(4).Range.Text = \"\u05db\u05d0\u05de\u05d5\u05e8 \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05dc\u05d4\u05db\u05e0\u05d9\u05e1 \u05d3\u05d9\u05e0\u05de\u05d9\u05d5\u05ea \u05d1\u05d1\u05e0\u05d9\u05d9\u05ea \u05e4\u05d0\u05e6'\u05d9\u05dd (\u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05d7\u05d3\u05e9\u05d9\u05dd). \u05d0\u05d9\u05da \u05d4\u05d5\u05d0 \u05e2\u05d5\u05e9\u05d4 \u05d0\u05ea \u05d6\u05d4. \u05dc\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e0\u05ea\u05d5\u05df \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc \u05e8\u05d3\u05d5\u05d3 \u05d9\u05d7\u05e1\u05d9\u05ea \u05d1\u05e8\u05de\u05d4 \u05e9\u05dc \u05d1\u05d8\u05d9\u05dd (bytes) \u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05d8\u05e8\u05d4 \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05d9\u05d0 \u05dc\u05d7\u05d6\u05d5\u05ea \u05d0\u05ea \u05d4\u05d1\u05d9\u05d9\u05d8 \u05d4\u05d1\u05d0. \u05d5\u05d0\u05d6 \u05d1\u05de\u05d5\u05d3\u05dc \u05d4\u05d2\u05d3\u05d5\u05dc \u05e9\u05dc\u05e0\u05d5 \u05d4\u05dd \u05e7\u05d5\u05d1\u05e2\u05d9\u05dd \u05d0\u05ea \u05d2\u05d1\u05d5\u05dc\u05d5\u05ea \u05d4\u05e4\u05d0\u05e5 \u05e2\u05dc \u05e1\u05de\u05da \u05d0\u05e0\u05d8\u05e8\u05d5\u05e4\u05d9\u05d4 \u05e9\u05dc \u05d4\u05d1\u05d8\u05d9\u05dd. \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05dd \u05d4\u05d0\u05e0\u05d8\u05e8\u05d5\u05e4\u05d9\u05d4 \u05e9\u05dc \u05d4\u05d1\u05d9\u05d9\u05d8 \u05d0\u05d5 \u05d2\u05d3\u05d5\u05dc\u05d4 \u05de\u05e1\u05e3 \u05de\u05e1\u05d5\u05d9\u05dd \u05d0\u05d5 \u05d7\u05d5\u05d5\u05ea\u05d4 \u05e2\u05dc\u05d9\u05d4 \u05de\u05e2\u05dc \u05e1\u05e3 \u05de\u05e1\u05d5\u05d9\u05dd \u05de\u05e2\u05dc \u05d4\u05d0\u05e0\u05d8\u05e8\u05d5\u05e4\u05d9\u05d4 \u05e9\u05dc \u05d4\u05d1\u05d9\u05d9\u05d8 \u05d4\u05d1\u05d0, \u05e4\u05d5\u05ea\u05d7\u05d9\u05dd \u05e4\u05d0\u05e5' \u05d7\u05d3\u05e9. \u05d0\u05d7\u05e8\u05ea \u05de\u05de\u05e9\u05d9\u05db\u05d9\u05dd \u05d0\u05ea \u05d4\u05e4\u05d0\u05e5' \u05d4\u05e0\u05d5\u05db\u05d7\u05d9.\"\n$d.Paragraphs.Item(5).Range.Text = \"\u05d0\u05d1\u05dc \u05d0\u05d9\u05da \u05db\u05dc \u05d4\u05e1\u05d9\u05e4\u05d5\u05e8 \u05d4\u05d6\u05d4 \u05e2\u05d5\u05d1\u05d3 - \u05db\u05de\u05d5 \u05e9\u05d0\u05de\u05e8\u05ea\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05d5\u05d0 byte-level \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05d5\u05d0 \u05de\u05d0\u05d5\u05de\u05df \u05dc\u05d7\u05d6\u05d5\u05ea \u05d0\u05ea \u05d4\u05d1\u05d9\u05d9\u05d8 \u05d4\u05d1\u05d0 \u05d1\u05d8\u05e7\u05e1\u05d8. \u05d0\u05d1\u05dc \u05d1\u05de\u05e7\u05d5\u05dd \u05dc\u05d4\u05e1\u05ea\u05db\u05dc \u05e2\u05dc \u05d4\u05e7\u05d5\u05e0\u05e7\u05e1\u05d8 \u05d1\u05ea\u05d5\u05e8 \u05de\u05e2\u05e8\u05da \u05e9\u05dc \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05dc\u05d4\u05d7\u05dc\u05d9\u05e3 \u05d0\u05d5\u05ea\u05d5 \u05d1\u05e4\u05d0\u05e6\u05d9\u05dd \u05d3\u05d9\u05e0\u05de\u05d9\u05d9\u05dd \u05e0\u05e7\u05d1\u05e2\u05d9\u05dd \u05e2\u05dc \u05e1\u05de\u05da \u05d4\u05d0\u05e0\u05d8\u05e8\u05d5\u05e4\u05d9\u05d4 \u05db\u05de\u05d5 \u05e9\u05d4\u05e1\u05d1\u05e8\u05ea\u05d9 \u05e7\u05d5\u05d3\u05dd. \"\n$d.Paragraphs.Item(6).Range.Text = \"\u05d1\u05e0\u05d5\u05e1\u05e3 \u05dc\u05e4\u05d0\u05e6\u05d9\u05dd \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e9\u05ea\u05de\u05e9 \u05d2\u05dd \u05d1\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05d1\u05d8\u05d9\u05dd \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea n-grams (\u05dc\u05d5\u05e7\u05d7\u05d9\u05dd n-grams \u05dc\u05d1\u05d9\u05d9\u05d8 \u05e0\u05ea\u05d5\u05df \u05de n=3 \u05e2\u05d3 n=8, \u05de\u05e4\u05e2\u05dc\u05d9\u05dd \u05d0\u05d9\u05d6\u05d4 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05d0\u05e9, \u05e1\u05d5\u05db\u05de\u05d9\u05dd \u05d5\u05de\u05e0\u05e8\u05de\u05dc\u05d9\u05dd). \u05d0\u05ea \u05d4\u05ea\u05d5\u05e6\u05d0\u05d4 \u05d4\u05d5\u05e4\u05db\u05d9\u05dd \u05dc\u05d5\u05d5\u05e7\u05d8\u05d5\u05e8 (\u05d4\u05de\u05d0\u05de\u05e8 \u05dc\u05d0 \u05de\u05e4\u05e8\u05e9 \u05d0\u05d9\u05da- \u05e8\u05e7 \u05de\u05d6\u05db\u05d9\u05e8 \u05e9\u05d9\u05e9 \u05d0\u05d9\u05d6\u05d5 \u05e9\u05db\u05d1\u05d4 \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea \u05d4\u05de\u05e2\u05d5\u05e8\u05d1\u05ea \u05d1\u05d6\u05d4) \u05d5\u05de\u05d6\u05d9\u05df \u05d0\u05d5\u05ea\u05d5 \u05dc\u05de\u05d4 \u05e9\u05e7\u05e8\u05d5\u05d9 \u05d1\u05de\u05d0\u05de\u05e8 Encoder Multi-Headed Cross-Attention (\u05e0\u05e7\u05e8\u05d0 \u05dc\u05d6\u05d4 \u05dc\u05e4\u05e9\u05d8\u05d5\u05ea EMHCA). \"\n$d.Paragraphs.Item(7).Range.Text = \"\u05de\u05d8\u05e8\u05ea\u05d5 \u05e9\u05dc EMHCA \u05d4\u05d9\u05d0 \u05dc\u05e9\u05dc\u05d1 \u05d0\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05e4\u05d0\u05e6'\u05d9\u05dd \u05e2\u05dd \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d1\u05d8\u05d9\u05dd \u05e9\u05dc\u05d4\u05dd(\u05db\u05dc \u05e4\u05d0\u05e5 \u05de\u05ea\u05d7\u05e9\u05d1 \u05e8\u05e7 \u05d1\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d1\u05d8\u05d9\u05dd \u05e9\u05dc\u05d5 \u05d5\u05dc\u05d0 \u05e9\u05dc \u05d4\u05d0\u05d7\u05e8\u05d9\u05dd). \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05d4\u05ea\u05d7\u05dc\u05ea\u05d9 \u05e9\u05dc \u05db\u05dc \u05e4\u05d0\u05e5 \u05de\u05d7\u05d5\u05e9\u05d1 \u05db-pooling (\u05db\u05dc\u05d5\u05de\u05e8 \u05de\u05de\u05d5\u05e6\u05e2) \u05e9\u05dc \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d1\u05d8\u05d9\u05dd \u05e9\u05dc\u05d5 (\u05e0\u05d6\u05db\u05d9\u05e8 \u05d6\u05d4 \u05db\u05dc \u05e4\u05d0\u05e5 \u05d4\u05d9\u05e0\u05d5 \u05de\u05e2\u05e8\u05da \u05e9\u05dc \u05d4\u05d1\u05d8\u05d9\u05dd). \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05e0\u05d5 \u05d1\u05d5\u05e0\u05d9\u05dd \u05db\u05db\u05d4 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05db\u05dc \u05e4\u05d0\u05e5' \u05d4\u05de\u05ea\u05d7\u05e9\u05d1 \u05e8\u05e7 \u05d1\u05de\u05d4 \u05e9\u05d9\u05e9 \u05d1\u05ea\u05d5\u05db\u05d5 (internal representation).  \"\n$d.Paragraphs.Item(8).Range.Text = \"\u05d0\u05d6 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05d1\u05d8\u05d9\u05dd \u05d5\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05e4\u05d0\u05e6'\u05d9\u05dd \u05de\u05d5\u05d6\u05e0\u05d9\u05dd \u05dc-EMHCA \u05e9\u05d6\u05d4 \u05dc\u05de\u05e2\u05e9\u05d4 \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05d3\u05d9 \u05e8\u05d3\u05d5\u05d3 (\u05e2\u05dd \u05de\u05e2\u05d8 \u05e9\u05db\u05d1\u05d5\u05ea) \u05e9\u05de\u05d8\u05e8\u05d5\u05ea \u05dc\u05d1\u05e0\u05d5\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2 \u05ea\u05dc\u05d5\u05d9 \u05d4\u05e7\u05e9\u05e8 \u05e9\u05e4\u05d0\u05e6'\u05d9\u05dd \u05db\u05ea\u05dc\u05d5\u05ea \u05d1\u05d1\u05d8\u05d9\u05dd \u05e9\u05dc\u05d5. \u05db\u05dc\u05d5\u05de\u05e8 \u05d2\u05dd \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d1\u05d8\u05d9\u05dd \u05d4\u05dd keys and values \u05db\u05d0\u05df \u05db\u05d0\u05e9\u05e8 \u05d4-queries \u05d4\u05dd \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05e4\u05d0\u05e6\u05d9\u05dd. \u05db\u05d0\u05de\u05d5\u05e8 \u05de\u05d4 \u05e9\u05d9\u05d5\u05e6\u05d0 \u05de\u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05d4\u05e8\u05d3\u05d5\u05d3 \u05d4\u05d6\u05d4 \u05d4\u05d5\u05d0 \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05e4\u05d0\u05e6'\u05d9\u05dd. \u05e0\u05e6\u05d9\u05d9\u05df \u05e9- EMHCA \u05e4\u05d5\u05dc\u05d8 \u05d2\u05dd \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d1\u05d9\u05d8\u05d9\u05dd \u05d1\u05e1\u05d5\u05e3 (\u05dc\u05d0 \u05d4\u05e6\u05dc\u05d7\u05ea\u05d9 \u05dc\u05d4\u05d1\u05d9\u05df \u05d0\u05d9\u05da \u05d6\u05d4 \u05e0\u05d1\u05e0\u05d4). \"\n$d.Paragraphs.Item(9).Range.Text = \"\u05db\u05dc \u05d0\u05dc\u05d5 \u05de\u05d5\u05db\u05e0\u05e1\u05d9\u05dd \u05dc\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05d9\u05d5\u05ea\u05e8 \u05e2\u05de\u05d5\u05e7 \u05d5\u05db\u05d1\u05d3 \u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05ea \u05d4\u05d9\u05d5\u05e6\u05e8 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d9\u05d5\u05ea\u05e8 \u05f4\u05e2\u05de\u05d5\u05e7\u05f4 \u05e9\u05dc \u05d4\u05e4\u05d0\u05e6\u05d9\u05dd. \u05d1\u05e9\u05dc\u05d1 \u05d4\u05d0\u05d7\u05e8\u05d5\u05df \u05d9\u05e9 \u05d0\u05ea \u05d4-Local Decoder \u05e9\u05d4\u05d5\u05e4\u05da \u05d0\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05e4\u05d0\u05e6'\u05d9\u05dd \u05d9\u05d7\u05d3 \u05e2\u05dd \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d1\u05d8\u05d9\u05dd \u05dc\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d1\u05d8\u05d9\u05dd \u05d4\u05e1\u05d5\u05e4\u05d9\u05d9\u05dd \u05e9\u05de\u05d4\u05dd \u05e0\u05d7\u05d6\u05d4 \u05d4\u05d1\u05d9\u05d9\u05d8 \u05d4\u05d1\u05d0. \u05d6\u05d4 \u05d2\u05dd \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05e8\u05d3\u05d5\u05d3 \u05d0\u05d1\u05dc \u05d4\u05e4\u05e2\u05dd \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05e4\u05d0\u05e6'\u05d9\u05dd \u05d4\u05dd keys and values \u05d5\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d1\u05d8\u05d9\u05dd \u05d4\u05dd \u05d4-queries.\"\n\n# --- Append two brand-new paragraphs, then a brand-new arXiv link paragraph,\n# after the paragraph that used to hold the (now-replaced) old link text. ---\n$d.Paragraphs.Item(9).Range.InsertParagraphAfter()\n$d.Paragraphs.Item(10).Range.Text = \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d8\u05d5\u05e2\u05df \u05dc\u05db\u05dc \u05de\u05d9\u05e0\u05d9 \u05d9\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05e9\u05dc \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05de\u05d5\u05e6\u05e2\u05ea \u05db\u05de\u05d5 \u05d9\u05db\u05d5\u05dc\u05ea \u05dc\u05d7\u05d6\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05dc\u05e2\u05dc\u05d5\u05ea \u05d0\u05d9\u05e0\u05e4\u05e8\u05e0\u05e1 \u05e7\u05d1\u05d5\u05e2\u05d4, \u05d5\u05de\u05e6\u05d9\u05d2\u05d4 \u05d3\u05d9\u05d5\u05e7 \u05de\u05e9\u05d5\u05e4\u05e8 \u05d1\u05d0\u05d9\u05de\u05d5\u05df \u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd.\"\n\n$d.Paragraphs.Item(10).Range.InsertParagraphAfter()\n$d.Paragraphs.Item(11).Range.Text = \"\u05d0\u05d5\u05e7\u05d9\u05d9, \u05d7\u05d9\u05d9\u05d1 \u05dc\u05d4\u05d2\u05d9\u05d3 \u05e9\u05d4\u05de\u05d0\u05de\u05e8 \u05dc\u05d0 \u05db\u05ea\u05d5\u05d1 \u05db\u05d6\u05d4 \u05d8\u05d5\u05d1 - \u05d9\u05e9 \u05d3\u05d1\u05e8\u05d9\u05dd \u05e9\u05dc\u05d0 \u05d4\u05d5\u05e1\u05d1\u05e8\u05d5 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d1\u05e8\u05d5\u05e8\u05d4 (\u05dc\u05de\u05d9\u05d8\u05d1 \u05d9\u05d3\u05d9\u05e2\u05ea\u05d9 \u05db\u05de\u05d5\u05d1\u05df). \u05d0\u05e0\u05d9 \u05e8\u05e7 \u05de\u05e7\u05d5\u05d5\u05d4 \u05e9\u05d4\u05e6\u05dc\u05d7\u05ea\u05d9 \u05dc\u05d4\u05d1\u05d9\u05df \u05d0\u05d5\u05ea\u05d5 \u05e0\u05db\u05d5\u05df\u2026.\"\n\n$d.Paragraphs.Item(11).Range.InsertParagraphAfter()\n$d.Paragraphs.Item(12).Range.Text = \"https://arxiv.org/abs/2412.09871\"\n"}
